$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 156
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 27
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 2
